$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-78 down to 18-79.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new data record.
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44707
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 100112035
$ws.Range("G17").Value = "Bruselas (repollito)"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 30000
$ws.Range("N17").Value = '$/malla 10 kilos'
$ws.Range("O17").Value = "Provincia de Quillota"
$ws.Range("P17").Value = 3000
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = "Hortaliza"
